# feat: neural networks with optimization of parameters
#
# Adds a second worksheet ("neural_networks") with the results of the
# neural-network hyper-parameter search, mirroring the formatting of the
# existing "Sheet1" (SVM results) worksheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Move the selection on the existing sheet first (selecting a range makes
# its sheet the active one, so this must happen before we add/activate the
# new sheet below).
$ws1.Range("A18").Select() | Out-Null

# Duplicate Sheet1 so the new sheet inherits the same namespaces / default
# row height / column-based formatting, then rename + wipe its contents.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "neural_networks"
$ws2.Cells.ClearContents()

# Column widths: col A = 29.5 chars, cols B:C = 19.83203125 chars (best-fit
# in the original). Excel pads the "ColumnWidth" you set by ~5/6 of a
# character when it stores the sheet's <col> width, so compensate here.
$pad = 0.8333333333333334
$ws2.Columns.Item(1).ColumnWidth = 29.5 - $pad
$ws2.Range("B1:C1").EntireColumn.ColumnWidth = 19.83203125 - $pad

$data = @(
    @("metoda",         "F1",                    "Accuracy"),
    @("adam-identity",  "0.64797507788161979",   "0.54800000000000004"),
    @("adam-logistic",  "0.77551020408163263",   "0.78000000000000003"),
    @("lbfgs-logistic", "0.69795918367346943",   "0.70399999999999996"),
    @("lbfgs-tanh",     "0.72332015810276684",   "0.71999999999999997"),
    @("adam-relu",      "0.71017274472168901",   "0.69799999999999995"),
    @("lbfgs-relu",     "0.63752276867030966",   "0.60199999999999998"),
    @("lbfgs-identity", "0.61056105610561062",   "0.52800000000000002"),
    @("adam-tanh",      "0.79918032786885251",   "0.80400000000000005"),
    @("sgd-identity",   "0.0",                   "0.496"),
    @("sgd-tanh",       "0.6851485148514852",    "0.68200000000000005"),
    @("sgd-relu",       "0.67021276595744683",   "0.504"),
    @("sgd-logistic",   "0.0",                   "0.496")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 1
    $row = $data[$i]

    $ws2.Cells.Item($r, 1).Value = $row[0]

    # Columns B/C hold numeric-looking values that must be stored as TEXT
    # (matching the source workbook). Force text via NumberFormat, then
    # reinstate the default "Normal" style so no stray quote-prefix style
    # is left attached to the cell.
    $cellB = $ws2.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $row[1]
    $cellB.Style = "Normal"

    $cellC = $ws2.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $row[2]
    $cellC.Style = "Normal"
}

# Select A3 on the new sheet to match the saved selection / make it the
# active (tab-selected) sheet.
$ws2.Range("A3").Select() | Out-Null
